# Insert a new weekly price observation row at row 319 of "Sheet1".
# This shifts all existing data rows 319..390 down to 320..391
# (dimension grows from A1:R390 to A1:R391), and fills the newly
# inserted row with the latest data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 319, pushing the rest down.
$ws.Rows.Item(319).Insert()

# Populate the freshly inserted row with the new record's data.
$ws.Cells.Item(319, 1).Value = 3
$ws.Cells.Item(319, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(319, 3).Value = "Coquimbo"
$ws.Cells.Item(319, 4).Value = 44785
$ws.Cells.Item(319, 5).Value = 5
$ws.Cells.Item(319, 6).Value = 100112009
$ws.Cells.Item(319, 7).Value = "Acelga"
$ws.Cells.Item(319, 8).Value = "Sin especificar"
$ws.Cells.Item(319, 9).Value = "Primera"
$ws.Cells.Item(319, 10).Value = 230
$ws.Cells.Item(319, 11).Value = 3000
$ws.Cells.Item(319, 12).Value = 3300
$ws.Cells.Item(319, 13).Value = 3157
$ws.Cells.Item(319, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(319, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(319, 16).Value = 526
$ws.Cells.Item(319, 17).Value = 6
$ws.Cells.Item(319, 18).Value = "Hortaliza"
